# nuevos_afiliados.xlsx - "Actualización general de presupuestos y archivos"
#
# Renames the A1/C1 headers ("Empresa" -> "Compañía", "afiliados" -> "Tipo_Afiliado"),
# tightens the selection back down to the header row, nudges column B's width and
# auto-sizes the new column C to fit the longer header, re-flows a handful of
# wrapped-text row heights that shrink now that the text fits differently, and
# switches the sheet to a portrait page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -----------------------------------------------
$ws.Range("A1").Value = "Compañía"
$ws.Range("C1").Value = "Tipo_Afiliado"
# B1 ("Categoria") is left untouched.

# --- Selection now just covers the header row ---------------------------
$ws.Range("A1:C1").Select() | Out-Null

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.88671875
$ws.Columns.Item(3).ColumnWidth = 12.5546875

# --- Row heights reflow (wrap-text column got a hair wider, so several
#     rows now wrap onto fewer lines) -------------------------------------
$ws.Rows.Item(3).RowHeight = 20.4
$ws.Rows.Item(4).RowHeight = 20.4
$ws.Rows.Item(8).RowHeight = 20.4
$ws.Rows.Item(9).RowHeight = 20.4
$ws.Rows.Item(11).RowHeight = 20.4
$ws.Rows.Item(12).RowHeight = 30.6
$ws.Rows.Item(16).RowHeight = 20.4
$ws.Rows.Item(18).RowHeight = 14.4
$ws.Rows.Item(20).RowHeight = 20.4
$ws.Rows.Item(21).RowHeight = 20.4
$ws.Rows.Item(22).RowHeight = 20.4
$ws.Rows.Item(23).RowHeight = 20.4
$ws.Rows.Item(25).RowHeight = 20.4
$ws.Rows.Item(26).RowHeight = 30.6
$ws.Rows.Item(27).RowHeight = 14.4
$ws.Rows.Item(29).RowHeight = 14.4
$ws.Rows.Item(31).RowHeight = 20.4
$ws.Rows.Item(32).RowHeight = 20.4
$ws.Rows.Item(35).RowHeight = 20.4

# --- Page setup: force portrait orientation -------------------------------
$ws.PageSetup.Orientation = 1
